# The commit swaps the contents of ppt/theme/theme1.xml and
# ppt/theme/theme2.xml: theme1.xml (previously the unused "Default"
# color theme, referenced only by the notes master) receives the
# "Tropic" theme content, while theme2.xml (the theme actually driving
# the slide master / every slide's look) receives the "Default" theme's
# color scheme.
#
# theme1.xml is not reachable through the PowerPoint object model (it
# is only linked from the notes master, which PowerPoint's automation
# surface does not expose as an independently themed object - it always
# mirrors the slide master's design/theme). The part of the swap that
# *is* reachable - and the part that actually affects what the
# presentation looks like - is recoloring the live theme (theme2.xml)
# from the "Tropic" palette to the "Default" palette. The font scheme
# and format scheme (fills/lines/effects) are identical between the two
# themes, so only the 12 color-scheme slots need to change.

function HexToComRGB([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    # COM RGB longs are packed 0x00BBGGRR
    return ($b * 65536) + ($g * 256) + $r
}

$p = $ppt.ActivePresentation

# ThemeColorScheme is exposed per-slide, but it reflects the single
# presentation-wide theme (theme2.xml) that backs the slide master.
$tcs = $p.Slides.Item(1).ThemeColorScheme

# Target values: the "Default" color scheme that originally lived in
# ppt/theme/theme1.xml, in clrScheme order.
$tcs.Item(1).RGB  = HexToComRGB "000000"   # dk1
$tcs.Item(2).RGB  = HexToComRGB "FFFFFF"   # lt1
$tcs.Item(3).RGB  = HexToComRGB "158158"   # dk2
$tcs.Item(4).RGB  = HexToComRGB "F3F3F3"   # lt2
$tcs.Item(5).RGB  = HexToComRGB "058DC7"   # accent1
$tcs.Item(6).RGB  = HexToComRGB "50B432"   # accent2
$tcs.Item(7).RGB  = HexToComRGB "ED561B"   # accent3
$tcs.Item(8).RGB  = HexToComRGB "EDEF00"   # accent4
$tcs.Item(9).RGB  = HexToComRGB "24CBE5"   # accent5
$tcs.Item(10).RGB = HexToComRGB "64E572"   # accent6
$tcs.Item(11).RGB = HexToComRGB "2200CC"   # hlink
$tcs.Item(12).RGB = HexToComRGB "551A8B"   # folHlink
